# -----------------------------------------------------------------------
# AutoLTS exp_records.xlsx -- add "Res50FC_speed" sheet, two new rows to
# Res50FC, and some new cell values / view state on the other sheets.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$wsRes50FC  = $wb.Worksheets.Item("Res50FC")
$wsMoCo     = $wb.Worksheets.Item("MoCo")
$wsMoCoClf  = $wb.Worksheets.Item("MoCoClf")

# -------------------------------------------------------------------
# 1) Res50FC: two new experiment rows (23, 24)
# -------------------------------------------------------------------
$wsRes50FC.Range("A20:E20").Copy() | Out-Null
$wsRes50FC.Range("A23:E23").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("G3:I3").Copy() | Out-Null
$wsRes50FC.Range("G23:I23").PasteSpecial(-4122) | Out-Null

$wsRes50FC.Range("A23").Value = 8857319
$wsRes50FC.Range("B23").Value = "Res50FC (for speed)"
$wsRes50FC.Range("C23").Value = 0.01
$wsRes50FC.Range("D23").Value = "SGD"
$wsRes50FC.Range("E23").Value = 128
$wsRes50FC.Range("G23").Value = $true
$wsRes50FC.Range("H23").Value = $true
$wsRes50FC.Range("I23").Value = $false

$wsRes50FC.Range("A20:E20").Copy() | Out-Null
$wsRes50FC.Range("A24:E24").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("G3:I3").Copy() | Out-Null
$wsRes50FC.Range("G24:I24").PasteSpecial(-4122) | Out-Null

$wsRes50FC.Range("A24").Value = 8857317
$wsRes50FC.Range("B24").Value = "Res50FC (for speed)"
$wsRes50FC.Range("C24").Value = 0.01
$wsRes50FC.Range("D24").Value = "SGD"
$wsRes50FC.Range("E24").Value = 128
$wsRes50FC.Range("G24").Value = $false
$wsRes50FC.Range("H24").Value = $true
$wsRes50FC.Range("I24").Value = $false

Write-Host "step1 done"

# -------------------------------------------------------------------
# 2) New sheet "Res50FC_speed" inserted between Res50FC and MoCo
# -------------------------------------------------------------------
$wsSpeed = $wb.Worksheets.Add($wsMoCo)
$wsSpeed.Name = "Res50FC_speed"

# --- header row (copy format + values from Res50FC row 1, cols A:L) ---
$wsRes50FC.Range("A1:L1").Copy() | Out-Null
$wsSpeed.Range("A1:L1").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("M1:N1").Copy() | Out-Null
$wsSpeed.Range("M1:N1").PasteSpecial(-4122) | Out-Null

$wsSpeed.Range("A1").Value = "slurm_job_id"
$wsSpeed.Range("B1").Value = "model"
$wsSpeed.Range("C1").Value = "learning rate"
$wsSpeed.Range("D1").Value = "optimizer"
$wsSpeed.Range("E1").Value = "batch_size"
$wsSpeed.Range("F1").Value = "best_epoch_id (starting from 0)"
$wsSpeed.Range("G1").Value = "frozen"
$wsSpeed.Range("H1").Value = "Augmentation"
$wsSpeed.Range("I1").Value = "Biased_sampling"
$wsSpeed.Range("J1").Value = "train_accuracy"
$wsSpeed.Range("K1").Value = "vali_accuracy"
$wsSpeed.Range("L1").Value = "test_accuracy"

Write-Host "step2 header done"

# --- data rows ---------------------------------------------------
# row 2 (full 14-col style A:N matches Res50FC row 20 pattern + J:L via
# a style-3 cell + M:N via a style-8 cell)
$wsRes50FC.Range("A20:F20").Copy() | Out-Null
$wsSpeed.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("G3:I3").Copy() | Out-Null
$wsSpeed.Range("G2:I2").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("A6:C6").Copy() | Out-Null
$wsSpeed.Range("J2:L2").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("M12:N12").Copy() | Out-Null
$wsSpeed.Range("M2:N2").PasteSpecial(-4122) | Out-Null

$wsSpeed.Range("A2").Value = 8857317
$wsSpeed.Range("B2").Value = "Res50FC (for speed)"
$wsSpeed.Range("C2").Value = 0.01
$wsSpeed.Range("D2").Value = "SGD"
$wsSpeed.Range("E2").Value = 128
$wsSpeed.Range("F2").ClearContents()
$wsSpeed.Range("G2").Value = $false
$wsSpeed.Range("H2").Value = $true
$wsSpeed.Range("I2").Value = $false
$wsSpeed.Range("J2").ClearContents()
$wsSpeed.Range("K2").ClearContents()
$wsSpeed.Range("L2").ClearContents()
$wsSpeed.Range("M2").Value = "fluctuating validation loss"
$wsSpeed.Range("N2").ClearContents()

# row 3: A3 has NO explicit style (plain, default)
$wsRes50FC.Range("B20:F20").Copy() | Out-Null
$wsSpeed.Range("B3:F3").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("G3:I3").Copy() | Out-Null
$wsSpeed.Range("G3:I3").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("M12").Copy() | Out-Null
$wsSpeed.Range("M3").PasteSpecial(-4122) | Out-Null

$wsSpeed.Range("A3").Value = 8857726
$wsSpeed.Range("B3").Value = "Res50FC (for speed)"
$wsSpeed.Range("C3").Value = 0.001
$wsSpeed.Range("D3").Value = "SGD"
$wsSpeed.Range("E3").Value = 128
$wsSpeed.Range("F3").ClearContents()
$wsSpeed.Range("G3").Value = $false
$wsSpeed.Range("H3").Value = $true
$wsSpeed.Range("I3").Value = $false
$wsSpeed.Range("M3").Value = "fluctuating validation loss"

# row 4: A4 has NO explicit style (plain, default)
$wsRes50FC.Range("B20:F20").Copy() | Out-Null
$wsSpeed.Range("B4:F4").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("G3:I3").Copy() | Out-Null
$wsSpeed.Range("G4:I4").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("M12").Copy() | Out-Null
$wsSpeed.Range("M4").PasteSpecial(-4122) | Out-Null

$wsSpeed.Range("A4").Value = 8857843
$wsSpeed.Range("B4").Value = "Res50FC (for speed)"
$wsSpeed.Range("C4").Value = 0.0001
$wsSpeed.Range("D4").Value = "SGD"
$wsSpeed.Range("E4").Value = 128
$wsSpeed.Range("F4").ClearContents()
$wsSpeed.Range("G4").Value = $false
$wsSpeed.Range("H4").Value = $true
$wsSpeed.Range("I4").Value = $false
$wsSpeed.Range("M4").Value = "better but still fluctuating"

# row 5: A5 has NO explicit style (plain, default)
$wsRes50FC.Range("B20:F20").Copy() | Out-Null
$wsSpeed.Range("B5:F5").PasteSpecial(-4122) | Out-Null
$wsRes50FC.Range("G3:I3").Copy() | Out-Null
$wsSpeed.Range("G5:I5").PasteSpecial(-4122) | Out-Null

$wsSpeed.Range("A5").Value = 8858089
$wsSpeed.Range("B5").Value = "Res50FC (for speed)"
$wsSpeed.Range("C5").Value = 0.00001
$wsSpeed.Range("D5").Value = "SGD"
$wsSpeed.Range("E5").Value = 128
$wsSpeed.Range("F5").ClearContents()
$wsSpeed.Range("G5").Value = $false
$wsSpeed.Range("H5").Value = $true
$wsSpeed.Range("I5").Value = $false

Write-Host "step2 data rows done"


